{"js": "// Rean\u00e1lise do texto do artefato \"WEG MOTOR SCAN\" (01. Declara\u00e7\u00e3o do escopo).\n// Faz 3 pequenas corre\u00e7\u00f5es de reda\u00e7\u00e3o, preservando a formata\u00e7\u00e3o original\n// (o texto editado herda a formata\u00e7\u00e3o do trecho pesquisado em cada `range`).\n\nconst body = context.document.body;\n\n// Pequena fun\u00e7\u00e3o auxiliar: localiza `searchText` (ocorr\u00eancia \u00fanica e exata)\n// e substitui pelo novo texto, mantendo a formata\u00e7\u00e3o do range encontrado.\nasync function replaceOnce(searchText, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Esperava 1 ocorr\u00eancia de \"${searchText}\", encontrei ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \"revis\u00f5es tem que ser feitas periodicamente, com isso precisamos\" ->\n//    \"revis\u00f5es precisam ser feitas periodicamente, e, com isso, precisamos\"\nawait replaceOnce(\n  \"tem que ser feitas periodicamente, com isso precisamos\",\n  \"precisam ser feitas periodicamente, e, com isso, precisamos\"\n);\n\n// 2) \"e se ele ainda funciona. Essa\" -> \"e se ele continua operando normalmente. Essa\"\nawait replaceOnce(\n  \"e se ele ainda funciona. Essa\",\n  \"e se ele continua operando normalmente. Essa\"\n);\n\n// 3) \"muitas vezes s\u00f3 sabermos quando um motor est\u00e1 ruim quando\" ->\n//    \"muitas vezes s\u00f3 descobriremos que um motor estava com algum defeito quando\"\nawait replaceOnce(\n  \"muitas vezes s\u00f3 sabermos quando um motor est\u00e1 ruim quando\",\n  \"muitas vezes s\u00f3 descobriremos que um motor estava com algum defeito quando\"\n);\n\n// 4) \"o que gera mais gastos.\" -> \"o que gera ainda mais gastos.\"\nawait replaceOnce(\n  \"o que gera mais gastos.\",\n  \"o que gera ainda mais gastos.\"\n);\n\n// 5) \"Comunica\u00e7\u00e3o por Bluetooth, Gateway, WiFi, etc\" -> adiciona um ponto final.\nawait replaceOnce(\n  \"Comunica\u00e7\u00e3o por Bluetooth, Gateway, WiFi, etc\",\n  \"Comunica\u00e7\u00e3o por Bluetooth, Gateway, WiFi, etc.\"\n);\n\nawait context.sync();\n", "ps1": "# Rean\u00e1lise do texto do artefato \"WEG MOTOR SCAN\" (01. Declara\u00e7\u00e3o do escopo).\n# Faz 3 pequenas corre\u00e7\u00f5es de reda\u00e7\u00e3o via Find & Replace, preservando a\n# formata\u00e7\u00e3o original de cada trecho (Find.Execute mant\u00e9m o rPr do texto\n# localizado ao substituir).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute(\n        $findText,    # FindText\n        $false,       # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        2             # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"Texto nao encontrado: $findText\"\n    }\n}\n\n# 1) \"revis\u00f5es tem que ser feitas periodicamente, com isso precisamos\" ->\n#    \"revis\u00f5es precisam ser feitas periodicamente, e, com isso, precisamos\"\nReplace-Text \"tem que ser feitas periodicamente, com isso precisamos\" \"precisam ser feitas periodicamente, e, com isso, precisamos\"\n\n# 2) \"e se ele ainda funciona. Essa\" -> \"e se ele continua operando normalmente. Essa\"\nReplace-Text \"e se ele ainda funciona. Essa\" \"e se ele continua operando normalmente. Essa\"\n\n# 3) \"muitas vezes s\u00f3 sabermos quando um motor est\u00e1 ruim quando\" ->\n#    \"muitas vezes s\u00f3 descobriremos que um motor estava com algum defeito quando\"\nReplace-Text \"muitas vezes s\u00f3 sabermos quando um motor est\u00e1 ruim quando\" \"muitas vezes s\u00f3 descobriremos que um motor estava com algum defeito quando\"\n\n# 4) \"o que gera mais gastos.\" -> \"o que gera ainda mais gastos.\"\nReplace-Text \"o que gera mais gastos.\" \"o que gera ainda mais gastos.\"\n\n# 5) \"Comunica\u00e7\u00e3o por Bluetooth, Gateway, WiFi, etc\" -> adiciona ponto final.\nReplace-Text \"Comunica\u00e7\u00e3o por Bluetooth, Gateway, WiFi, etc\" \"Comunica\u00e7\u00e3o por Bluetooth, Gateway, WiFi, etc.\"\n"}
